# Re-run of the "resolve and classify+summarise" pipeline after a mapping-file
# change produced zeroed-out / recomputed summary numbers on several sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "Range Status": species counts collapse to 0 and the percentage
#     column is dropped entirely for rows 2-7 (Historical .. Very Large) ---
$wsRange = $wb.Worksheets.Item("Range Status")
for ($r = 2; $r -le 7; $r++) {
    $wsRange.Cells.Item($r, 2).Value = 0       # column B -> 0
    $wsRange.Cells.Item($r, 3).ClearContents()  # column C removed entirely
}

# --- Sheet "Species qualification": Range Analysis row species count -> 0 ---
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Cells.Item(5, 2).Value = 0   # B5 (Range Analysis) 452 -> 0

# --- Sheet "High Priority break-up": recomputed figures; row 3 ("Range")
#     becomes "IUCN" with new numbers, and the old row 4 ("IUCN") is removed ---
$wsBreak = $wb.Worksheets.Item("High Priority break-up")
$wsBreak.Cells.Item(2, 5).Value = 22.2   # E2 33.3 -> 22.2

$wsBreak.Cells.Item(3, 1).Value = "IUCN"
$wsBreak.Cells.Item(3, 2).Value = 7
$wsBreak.Cells.Item(3, 3).Value = 77.8
$wsBreak.Cells.Item(3, 4).Value = 7
$wsBreak.Cells.Item(3, 5).Value = 77.8

$wsBreak.Rows.Item(4).Delete()
